# "Mudança na SUPORTE 2": a client ("ANTONIO LUCIANO DE CAMARGO FILHO &
# MARCIA MIYUKI IOSHIHARA" / grupo "ANTONIO LUCIANO") was previously only
# listed on the onshore-only tables (clients_onshore / clients_onshore (2))
# with a single row. The edit removes that single-row entry from the
# onshore-only tables and instead lists the client properly (Onshore +
# Offshore pair) on the on_off tables (clients_on_off / clients_on_off (2)).

$wb = $excel.ActiveWorkbook

# --- clients_onshore: remove the old "ANTONIO LUCIANO" row (row 9) ---
$wsOnshore1 = $wb.Worksheets.Item("clients_onshore")
$wsOnshore1.Rows.Item(9).Delete()

# --- clients_onshore (2): mirror table, same removal ---
$wsOnshore2 = $wb.Worksheets.Item("clients_onshore (2)")
$wsOnshore2.Rows.Item(9).Delete()

# --- clients_on_off: fill the two already-reserved blank rows (76/77)
# with the Onshore/Offshore pair for this client ---
$wsOnOff1 = $wb.Worksheets.Item("clients_on_off")

$wsOnOff1.Range("A76").Value = "ANTONIO LUCIANO"
$wsOnOff1.Range("B76").Value = "ANTONIO LUCIANO DE CAMARGO FILHO & MARCIA MIYUKI IOSHIHARA"
$wsOnOff1.Range("C76").Value = 1010
$wsOnOff1.Range("D76").Value = "Onshore"
$wsOnOff1.Range("E76").Value = (Get-Date -Year 2024 -Month 3 -Day 22 -Hour 0 -Minute 0 -Second 0)

$wsOnOff1.Range("A77").Value = "ANTONIO LUCIANO"
$wsOnOff1.Range("B77").Value = "ANTONIO LUCIANO DE CAMARGO FILHO & MARCIA MIYUKI IOSHIHARA"
$wsOnOff1.Range("C77").Value = 2037
$wsOnOff1.Range("D77").Value = "Offshore"
$wsOnOff1.Range("E77").Value = (Get-Date -Year 2024 -Month 8 -Day 5 -Hour 0 -Minute 0 -Second 0)

# --- clients_on_off (2): mirror table; rows 76/77 don't exist yet, so
# first extend the table by copying the formatting from the last row
# (75) down, then set the same values ---
$wsOnOff2 = $wb.Worksheets.Item("clients_on_off (2)")

$wsOnOff2.Range("A75:E75").Copy()
$wsOnOff2.Range("A76:E77").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsOnOff2.Range("A76").Value = "ANTONIO LUCIANO"
$wsOnOff2.Range("B76").Value = "ANTONIO LUCIANO DE CAMARGO FILHO & MARCIA MIYUKI IOSHIHARA"
$wsOnOff2.Range("C76").Value = 1010
$wsOnOff2.Range("D76").Value = "Onshore"
$wsOnOff2.Range("E76").Value = (Get-Date -Year 2024 -Month 3 -Day 22 -Hour 0 -Minute 0 -Second 0)

$wsOnOff2.Range("A77").Value = "ANTONIO LUCIANO"
$wsOnOff2.Range("B77").Value = "ANTONIO LUCIANO DE CAMARGO FILHO & MARCIA MIYUKI IOSHIHARA"
$wsOnOff2.Range("C77").Value = 2037
$wsOnOff2.Range("D77").Value = "Offshore"
$wsOnOff2.Range("E77").Value = (Get-Date -Year 2024 -Month 8 -Day 5 -Hour 0 -Minute 0 -Second 0)

# --- the workbook's active tab moved from "clients_onshore (2)" to
# "clients_onshore" ---
$wsOnshore1.Activate()
